# Update crypto price/volume table: columns D (Price) and E (Volume 1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D updates whose new text is not a valid numeric literal (safe to set directly) ---
$ws.Range('D2').Value = '43.012.33'
$ws.Range('D3').Value = '2.310.97'
$ws.Range('D15').Value = '2.687.49'
$ws.Range('D16').Value = '2.348.91'
$ws.Range('D18').Value = '42.974.75'
$ws.Range('D21').Value = '0.0₃0904'
$ws.Range('D42').Value = '2.002.90'
$ws.Range('D49').Value = '2.530.85'

# --- Column D updates whose new text WOULD be auto-converted to a number by Excel. ---
# --- Force them to remain Text by applying a Text number format before assignment, ---
# --- then restore the default "Normal" style so no stray style index is left on the cell. ---
$textForcedCells = @('D5', 'D6', 'D10', 'D13', 'D14', 'D17', 'D19', 'D22', 'D23', 'D24', 'D27', 'D29', 'D30', 'D34', 'D36', 'D37', 'D43', 'D44', 'D45', 'D46', 'D48', 'D50')
foreach ($addr in $textForcedCells) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range('D5').Value = '303.32'
$ws.Range('D6').Value = '100.96'
$ws.Range('D10').Value = '35.25'
$ws.Range('D13').Value = '17.86'
$ws.Range('D14').Value = '6.92'
$ws.Range('D17').Value = '0.813'
$ws.Range('D19').Value = '12.68'
$ws.Range('D22').Value = '67.88'
$ws.Range('D23').Value = '237.33'
$ws.Range('D24').Value = '2.22'
$ws.Range('D27').Value = '24.77'
$ws.Range('D29').Value = '167.80'
$ws.Range('D30').Value = '34.04'
$ws.Range('D34').Value = '4.61'
$ws.Range('D36').Value = '17.00'
$ws.Range('D37').Value = '0.0692'
$ws.Range('D43').Value = '2.30'
$ws.Range('D44').Value = '0.0287'
$ws.Range('D45').Value = '10.23'
$ws.Range('D46').Value = '17.48'
$ws.Range('D48').Value = '55.43'
$ws.Range('D50').Value = '1.53'
foreach ($addr in $textForcedCells) { $ws.Range($addr).Style = "Normal" }

# --- Column E updates (percentage-change labels padded with spaces, always text) ---
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('E3').Value = '  +2.01%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('E6').Value = '  +5.88%  '
$ws.Range('E7').Value = '  +1.71%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +4.62%  '
$ws.Range('E10').Value = '  +5.89%  '
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('E12').Value = '  +4.14%  '
$ws.Range('E13').Value = '  +15.40%  '
$ws.Range('E14').Value = '  +4.02%  '
$ws.Range('E15').Value = '  +2.61%  '
$ws.Range('E16').Value = '  +3.23%  '
$ws.Range('E17').Value = '  +4.05%  '
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('E19').Value = '  +8.58%  '
$ws.Range('E20').Value = '  +3.35%  '
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('E22').Value = '  +2.11%  '
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('E24').Value = '  +13.57%  '
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  +3.44%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +2.58%  '
$ws.Range('E34').Value = '  +3.30%  '
$ws.Range('E35').Value = '  +3.68%  '
$ws.Range('E36').Value = '  +3.10%  '
$ws.Range('E37').Value = '  +1.13%  '
$ws.Range('E38').Value = '  +4.17%  '
$ws.Range('E39').Value = '  +4.42%  '
$ws.Range('E40').Value = '  +1.94%  '
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('E42').Value = '  +2.27%  '
$ws.Range('E43').Value = '  -5.63%  '
$ws.Range('E44').Value = '  +3.97%  '
$ws.Range('E45').Value = '  +7.45%  '
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('E47').Value = '  +2.73%  '
$ws.Range('E48').Value = '  +6.28%  '
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('E50').Value = '  +5.08%  '
$ws.Range('E51').Value = '  +0.99%  '

